$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the "Conversión del día" text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Value2
$text = $text.Replace("1000 Bs = 1.72 = 6300.79 pesos", "1000 Bs = 1.72 = 6312.52 pesos")
$text = $text.Replace("6300.79 pesos = 1.72 = 951.13 Bs", "6312.52 pesos = 1.72 = 955.92 Bs")
$cellA1.Value2 = $text

# --- Sheet "tasas": update the rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value2 = 580.6
$wsTasas.Range("O10").Value2 = 3665.05
$wsTasas.Range("N12").Value2 = 3665
$wsTasas.Range("O12").Value2 = 555
